# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1, matching the existing
# header formatting (bold font, thin border, centered) by copying the
# format from the last existing header cell (AC1).
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-47): every player row gets the team's season record.
$ws.Range("AD2:AD47").Value = 79
$ws.Range("AE2:AE47").Value = 83
$ws.Range("AF2:AF47").Value = 0
